$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = 3.7
$ws.Range("I11").Value = 1.8
$ws.Range("J11").Value = 3.95
$ws.Range("K11").Value = 2.45
$ws.Range("L11").Value = 2.2
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 9
$ws.Range("U11").Value = 1.65
$ws.Range("V11").Value = 1.98
$ws.Range("Y11").Value = 13
$ws.Range("AB11").Value = 37
$ws.Range("AD11").Value = 7.3
$ws.Range("AE11").Value = 14
$ws.Range("AF11").Value = 60
$ws.Range("AG11").Value = 8
$ws.Range("AH11").Value = 9.25
$ws.Range("AJ11").Value = 15
$ws.Range("AP11").Value = 21
$ws.Range("AS11").Value = 200
$ws.Range("AT11").Value = 3.6
$ws.Range("AU11").Value = 6.6
$ws.Range("AW11").Value = 4
$ws.Range("AZ11").Value = 24
